$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "fixing 3000 files from BES": the Q value was mistakenly entered as 10.6 and
# the observable label was missing the "-INT" suffix.

# 1) Rename the observable label "AUC-0-PT" -> "AUC-0-PT-INT" (column H, rows 2-7)
for ($r = 2; $r -le 7; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    if ($cell.Value2 -eq "AUC-0-PT") {
        $cell.Value = "AUC-0-PT-INT"
    }
}

# 2) Correct the Q value 10.6 -> 3.65 (column A, rows 2-7)
for ($r = 2; $r -le 7; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq 10.6) {
        $cell.Value = 3.65
    }
}
